$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (font/border/alignment) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# I0 / IF data values for rows 2-76
$ijValues = @{
  2 = @(4,4)
  3 = @(6,7)
  4 = @(5,6)
  5 = @(1,4)
  6 = @(7,7)
  7 = @(7,7)
  8 = @(5,6)
  9 = @(8,8)
  10 = @(7,7)
  11 = @(7,7)
  12 = @(5,6)
  13 = @(7,8)
  14 = @(10,10)
  15 = @(1,2)
  16 = @(6,6)
  17 = @(1,2)
  18 = @(9,9)
  19 = @(7,7)
  20 = @(1,2)
  21 = @(4,4)
  22 = @(6,6)
  23 = @(4,4)
  24 = @(9,9)
  25 = @(5,6)
  26 = @(5,6)
  27 = @(6,6)
  28 = @(5,5)
  29 = @(6,7)
  30 = @(6,6)
  31 = @(6,6)
  32 = @(6,6)
  33 = @(6,6)
  34 = @(9,9)
  35 = @(1,2)
  36 = @(8,8)
  37 = @(7,7)
  38 = @(8,8)
  39 = @(7,7)
  40 = @(7,7)
  41 = @(7,7)
  42 = @(6,6)
  43 = @(6,7)
  44 = @(7,8)
  45 = @(7,7)
  46 = @(6,7)
  47 = @(6,6)
  48 = @(7,7)
  49 = @(7,8)
  50 = @(5,5)
  51 = @(8,8)
  52 = @(7,7)
  53 = @(6,7)
  54 = @(6,7)
  55 = @(6,6)
  56 = @(6,6)
  57 = @(4,5)
  58 = @(6,6)
  59 = @(6,6)
  60 = @(4,5)
  61 = @(3,4)
  62 = @(9,9)
  63 = @(7,7)
  64 = @(7,9)
  65 = @(1,1)
  66 = @(1,4)
  67 = @(6,7)
  68 = @(5,7)
  69 = @(6,6)
  70 = @(7,7)
  71 = @(6,7)
  72 = @(3,4)
  73 = @(5,5)
  74 = @(4,4)
  75 = @(4,4)
  76 = @(3,3)
}

foreach ($r in $ijValues.Keys) {
    $pair = $ijValues[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
